# [Prototype] Add ex_PersonTimeline module
# Translate the PeopleEvents sheet (table TableEvents) from Russian labels
# to the English labels used by the new ex_PersonTimeline prototype module.
# Only the textual header/name/event-type values change - dates, numeric
# department/position/salary codes, and RecordNo all stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:F1). G1 "RecordNo" is already in English and unchanged. ---
$ws.Range("A1").Value = "FIO"
$ws.Range("B1").Value = "EventDate"
$ws.Range("C1").Value = "EventType"
$ws.Range("D1").Value = "Department"
$ws.Range("E1").Value = "Position"
$ws.Range("F1").Value = "Salary"

# --- Data rows 2-16: column A (full name) and column C (event type). ---
# Columns B (date), D (department), E (position), F (salary) and G (RecordNo)
# are untouched.
$names = @(
    "Ivanov Ivan Ivanovich",
    "Petrov Pyotr Petrovich",
    "Sidorova Anna Sergeevna",
    "Ivanov Ivan Ivanovich",
    "Shevchenko Oleg Andreevich",
    "Petrov Pyotr Petrovich",
    "Novikova Elena Viktorovna",
    "Ivanov Ivan Ivanovich",
    "Sidorova Anna Sergeevna",
    "Shevchenko Oleg Andreevich",
    "Petrov Pyotr Petrovich",
    "Ivanov Ivan Ivanovich",
    "Sidorova Anna Sergeevna",
    "Novikova Elena Viktorovna",
    "Shevchenko Oleg Andreevich"
)

$eventTypes = @(
    "Promotion",
    "Hired",
    "Hired",
    "Promotion",
    "Hired",
    "Transfer",
    "Hired",
    "Hired",
    "Transfer",
    "Transfer",
    "Promotion",
    "Transfer",
    "Promotion",
    "Promotion",
    "Promotion"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $eventTypes[$i]
}

# --- Selection moves to I10 (matches the saved cursor position in the diff). ---
[void]$ws.Range("I10").Select()
